$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 22.14384766666667
$ws.Cells.Item(2, 8).Value = 66.431543
$ws.Cells.Item(2, 9).Value = 0.05562336639723622
$ws.Cells.Item(2, 10).Value = 0.0556233663972362
$ws.Cells.Item(2, 13).Value = 2.655702666666667
$ws.Cells.Item(2, 14).Value = 7.967108
$ws.Cells.Item(2, 15).Value = 0.05887637219457465
$ws.Cells.Item(2, 16).Value = 0.05887637219457464
$ws.Cells.Item(2, 17).Value = 58.80747529862712
$ws.Cells.Item(2, 18).Value = 529.267277687644
$ws.Cells.Item(2, 19).Value = 0.003274902022718876
$ws.Cells.Item(2, 20).Value = 0.003274902022718875
$ws.Cells.Item(3, 7).Value = 22.14384766666667
$ws.Cells.Item(3, 8).Value = 66.431543
$ws.Cells.Item(3, 9).Value = 0.05562336639723622
$ws.Cells.Item(3, 10).Value = 0.0556233663972362
$ws.Cells.Item(3, 15).Value = 0.111327724745791
$ws.Cells.Item(3, 16).Value = 0.111327724745791
$ws.Cells.Item(3, 17).Value = 111.1974494862601
$ws.Cells.Item(3, 18).Value = 1000.777045376341
$ws.Cells.Item(3, 19).Value = 0.006192422823705796
$ws.Cells.Item(3, 20).Value = 0.006192422823705792
$ws.Cells.Item(4, 7).Value = 22.14384766666667
$ws.Cells.Item(4, 8).Value = 66.431543
$ws.Cells.Item(4, 9).Value = 0.05562336639723622
$ws.Cells.Item(4, 10).Value = 0.0556233663972362
$ws.Cells.Item(4, 13).Value = 0.9216453333333332
$ws.Cells.Item(4, 14).Value = 2.764936
$ws.Cells.Item(4, 15).Value = 0.02043268410948847
$ws.Cells.Item(4, 16).Value = 0.02043268410948846
$ws.Cells.Item(4, 17).Value = 20.40877386402756
$ws.Cells.Item(4, 18).Value = 183.678964776248
$ws.Cells.Item(4, 19).Value = 0.001136534674701063
$ws.Cells.Item(4, 20).Value = 0.001136534674701063
$ws.Cells.Item(5, 7).Value = 22.14384766666667
$ws.Cells.Item(5, 8).Value = 66.431543
$ws.Cells.Item(5, 9).Value = 0.05562336639723622
$ws.Cells.Item(5, 10).Value = 0.0556233663972362
$ws.Cells.Item(5, 13).Value = 35.230657
$ws.Cells.Item(5, 14).Value = 105.691971
$ws.Cells.Item(5, 15).Value = 0.7810562907612387
$ws.Cells.Item(5, 16).Value = 0.7810562907612385
$ws.Cells.Item(5, 17).Value = 780.1423018045838
$ws.Cells.Item(5, 18).Value = 7021.280716241253
$ws.Cells.Item(5, 19).Value = 0.04344498023787865
$ws.Cells.Item(5, 20).Value = 0.04344498023787862
$ws.Cells.Item(6, 7).Value = 22.14384766666667
$ws.Cells.Item(6, 8).Value = 66.431543
$ws.Cells.Item(6, 9).Value = 0.05562336639723622
$ws.Cells.Item(6, 10).Value = 0.0556233663972362
$ws.Cells.Item(6, 13).Value = 1.276824333333333
$ws.Cells.Item(6, 14).Value = 3.830473
$ws.Cells.Item(6, 15).Value = 0.02830692818890731
$ws.Cells.Item(6, 16).Value = 0.02830692818890731
$ws.Cells.Item(6, 17).Value = 28.27380353442656
$ws.Cells.Item(6, 18).Value = 254.464231809839
$ws.Cells.Item(6, 19).Value = 0.001574526638231846
$ws.Cells.Item(6, 20).Value = 0.001574526638231845
$ws.Cells.Item(7, 9).Value = 0.8709978578802913
$ws.Cells.Item(7, 10).Value = 0.8709978578802913
$ws.Cells.Item(7, 13).Value = 2.655702666666667
$ws.Cells.Item(7, 14).Value = 7.967108
$ws.Cells.Item(7, 15).Value = 0.05887637219457465
$ws.Cells.Item(7, 16).Value = 0.05887637219457464
$ws.Cells.Item(7, 17).Value = 920.8573362254718
$ws.Cells.Item(7, 18).Value = 8287.716026029248
$ws.Cells.Item(7, 19).Value = 0.05128119406123727
$ws.Cells.Item(7, 20).Value = 0.05128119406123726
$ws.Cells.Item(8, 9).Value = 0.8709978578802913
$ws.Cells.Item(8, 10).Value = 0.8709978578802913
$ws.Cells.Item(8, 15).Value = 0.111327724745791
$ws.Cells.Item(8, 16).Value = 0.111327724745791
$ws.Cells.Item(8, 19).Value = 0.09696620977627068
$ws.Cells.Item(8, 20).Value = 0.09696620977627066
$ws.Cells.Item(9, 9).Value = 0.8709978578802913
$ws.Cells.Item(9, 10).Value = 0.8709978578802913
$ws.Cells.Item(9, 13).Value = 0.9216453333333332
$ws.Cells.Item(9, 14).Value = 2.764936
$ws.Cells.Item(9, 15).Value = 0.02043268410948847
$ws.Cells.Item(9, 16).Value = 0.02043268410948846
$ws.Cells.Item(9, 17).Value = 319.5778944874239
$ws.Cells.Item(9, 18).Value = 2876.201050386815
$ws.Cells.Item(9, 19).Value = 0.01779682409010912
$ws.Cells.Item(9, 20).Value = 0.01779682409010912
$ws.Cells.Item(10, 9).Value = 0.8709978578802913
$ws.Cells.Item(10, 10).Value = 0.8709978578802913
$ws.Cells.Item(10, 13).Value = 35.230657
$ws.Cells.Item(10, 14).Value = 105.691971
$ws.Cells.Item(10, 15).Value = 0.7810562907612387
$ws.Cells.Item(10, 16).Value = 0.7810562907612385
$ws.Cells.Item(10, 17).Value = 12216.12997783886
$ws.Cells.Item(10, 18).Value = 109945.1698005498
$ws.Cells.Item(10, 19).Value = 0.6802983561369649
$ws.Cells.Item(10, 20).Value = 0.6802983561369647
$ws.Cells.Item(11, 9).Value = 0.8709978578802913
$ws.Cells.Item(11, 10).Value = 0.8709978578802913
$ws.Cells.Item(11, 13).Value = 1.276824333333333
$ws.Cells.Item(11, 14).Value = 3.830473
$ws.Cells.Item(11, 15).Value = 0.02830692818890731
$ws.Cells.Item(11, 16).Value = 0.02830692818890731
$ws.Cells.Item(11, 17).Value = 442.735201187632
$ws.Cells.Item(11, 18).Value = 3984.616810688688
$ws.Cells.Item(11, 19).Value = 0.0246552738157095
$ws.Cells.Item(11, 20).Value = 0.0246552738157095
$ws.Cells.Item(12, 7).Value = 0.1541363333333333
$ws.Cells.Item(12, 8).Value = 0.462409
$ws.Cells.Item(12, 9).Value = 0.0003871766945467397
$ws.Cells.Item(12, 10).Value = 0.0003871766945467395
$ws.Cells.Item(12, 13).Value = 2.655702666666667
$ws.Cells.Item(12, 14).Value = 7.967108
$ws.Cells.Item(12, 15).Value = 0.05887637219457465
$ws.Cells.Item(12, 16).Value = 0.05887637219457464
$ws.Cells.Item(12, 17).Value = 0.4093402714635556
$ws.Cells.Item(12, 18).Value = 3.684062443172
$ws.Cells.Item(12, 19).Value = 0.00002279555917319899
$ws.Cells.Item(12, 20).Value = 0.00002279555917319898
$ws.Cells.Item(13, 7).Value = 0.1541363333333333
$ws.Cells.Item(13, 8).Value = 0.462409
$ws.Cells.Item(13, 9).Value = 0.0003871766945467397
$ws.Cells.Item(13, 10).Value = 0.0003871766945467395
$ws.Cells.Item(13, 15).Value = 0.111327724745791
$ws.Cells.Item(13, 16).Value = 0.111327724745791
$ws.Cells.Item(13, 17).Value = 0.7740103435425556
$ws.Cells.Item(13, 18).Value = 6.966093091883
$ws.Cells.Item(13, 19).Value = 0.00004310350047848464
$ws.Cells.Item(13, 20).Value = 0.00004310350047848462
$ws.Cells.Item(14, 7).Value = 0.1541363333333333
$ws.Cells.Item(14, 8).Value = 0.462409
$ws.Cells.Item(14, 9).Value = 0.0003871766945467397
$ws.Cells.Item(14, 10).Value = 0.0003871766945467395
$ws.Cells.Item(14, 13).Value = 0.9216453333333332
$ws.Cells.Item(14, 14).Value = 2.764936
$ws.Cells.Item(14, 15).Value = 0.02043268410948847
$ws.Cells.Item(14, 16).Value = 0.02043268410948846
$ws.Cells.Item(14, 17).Value = 0.1420590323137778
$ws.Cells.Item(14, 18).Value = 1.278531290824
$ws.Cells.Item(14, 19).Value = 0.000007911059094229438
$ws.Cells.Item(14, 20).Value = 0.000007911059094229433
$ws.Cells.Item(15, 7).Value = 0.1541363333333333
$ws.Cells.Item(15, 8).Value = 0.462409
$ws.Cells.Item(15, 9).Value = 0.0003871766945467397
$ws.Cells.Item(15, 10).Value = 0.0003871766945467395
$ws.Cells.Item(15, 13).Value = 35.230657
$ws.Cells.Item(15, 14).Value = 105.691971
$ws.Cells.Item(15, 15).Value = 0.7810562907612387
$ws.Cells.Item(15, 16).Value = 0.7810562907612385
$ws.Cells.Item(15, 17).Value = 5.430324290904334
$ws.Cells.Item(15, 18).Value = 48.872918618139
$ws.Cells.Item(15, 19).Value = 0.0003024067929118736
$ws.Cells.Item(15, 20).Value = 0.0003024067929118734
$ws.Cells.Item(16, 7).Value = 0.1541363333333333
$ws.Cells.Item(16, 8).Value = 0.462409
$ws.Cells.Item(16, 9).Value = 0.0003871766945467397
$ws.Cells.Item(16, 10).Value = 0.0003871766945467395
$ws.Cells.Item(16, 13).Value = 1.276824333333333
$ws.Cells.Item(16, 14).Value = 3.830473
$ws.Cells.Item(16, 15).Value = 0.02830692818890731
$ws.Cells.Item(16, 16).Value = 0.02830692818890731
$ws.Cells.Item(16, 17).Value = 0.1968050210507778
$ws.Cells.Item(16, 18).Value = 1.771245189457
$ws.Cells.Item(16, 19).Value = 0.00001095978288895306
$ws.Cells.Item(16, 20).Value = 0.00001095978288895306
$ws.Cells.Item(17, 7).Value = 28.90575466666667
$ws.Cells.Item(17, 8).Value = 86.717264
$ws.Cells.Item(17, 9).Value = 0.07260867248616912
$ws.Cells.Item(17, 10).Value = 0.07260867248616912
$ws.Cells.Item(17, 13).Value = 2.655702666666667
$ws.Cells.Item(17, 14).Value = 7.967108
$ws.Cells.Item(17, 15).Value = 0.05887637219457465
$ws.Cells.Item(17, 16).Value = 0.05887637219457464
$ws.Cells.Item(17, 17).Value = 76.7650897502791
$ws.Cells.Item(17, 18).Value = 690.8858077525119
$ws.Cells.Item(17, 19).Value = 0.004274935225849665
$ws.Cells.Item(17, 20).Value = 0.004274935225849664
$ws.Cells.Item(18, 7).Value = 28.90575466666667
$ws.Cells.Item(18, 8).Value = 86.717264
$ws.Cells.Item(18, 9).Value = 0.07260867248616912
$ws.Cells.Item(18, 10).Value = 0.07260867248616912
$ws.Cells.Item(18, 15).Value = 0.111327724745791
$ws.Cells.Item(18, 16).Value = 0.111327724745791
$ws.Cells.Item(18, 17).Value = 145.1530123758631
$ws.Cells.Item(18, 18).Value = 1306.377111382768
$ws.Cells.Item(18, 19).Value = 0.008083358304697526
$ws.Cells.Item(18, 20).Value = 0.008083358304697525
$ws.Cells.Item(19, 7).Value = 28.90575466666667
$ws.Cells.Item(19, 8).Value = 86.717264
$ws.Cells.Item(19, 9).Value = 0.07260867248616912
$ws.Cells.Item(19, 10).Value = 0.07260867248616912
$ws.Cells.Item(19, 13).Value = 0.9216453333333332
$ws.Cells.Item(19, 14).Value = 2.764936
$ws.Cells.Item(19, 15).Value = 0.02043268410948847
$ws.Cells.Item(19, 16).Value = 0.02043268410948846
$ws.Cells.Item(19, 17).Value = 26.64085389501155
$ws.Cells.Item(19, 18).Value = 239.767685055104
$ws.Cells.Item(19, 19).Value = 0.0014835900685192
$ws.Cells.Item(19, 20).Value = 0.0014835900685192
$ws.Cells.Item(20, 7).Value = 28.90575466666667
$ws.Cells.Item(20, 8).Value = 86.717264
$ws.Cells.Item(20, 9).Value = 0.07260867248616912
$ws.Cells.Item(20, 10).Value = 0.07260867248616912
$ws.Cells.Item(20, 13).Value = 35.230657
$ws.Cells.Item(20, 14).Value = 105.691971
$ws.Cells.Item(20, 15).Value = 0.7810562907612387
$ws.Cells.Item(20, 16).Value = 0.7810562907612385
$ws.Cells.Item(20, 17).Value = 1018.368727987483
$ws.Cells.Item(20, 18).Value = 9165.318551887343
$ws.Cells.Item(20, 19).Value = 0.05671146040914486
$ws.Cells.Item(20, 20).Value = 0.05671146040914485
$ws.Cells.Item(21, 7).Value = 28.90575466666667
$ws.Cells.Item(21, 8).Value = 86.717264
$ws.Cells.Item(21, 9).Value = 0.07260867248616912
$ws.Cells.Item(21, 10).Value = 0.07260867248616912
$ws.Cells.Item(21, 13).Value = 1.276824333333333
$ws.Cells.Item(21, 14).Value = 3.830473
$ws.Cells.Item(21, 15).Value = 0.02830692818890731
$ws.Cells.Item(21, 16).Value = 0.02830692818890731
$ws.Cells.Item(21, 17).Value = 36.90757093176356
$ws.Cells.Item(21, 18).Value = 332.168138385872
$ws.Cells.Item(21, 19).Value = 0.002055328477957879
$ws.Cells.Item(21, 20).Value = 0.002055328477957879
$ws.Cells.Item(22, 5).Value = 3
$ws.Cells.Item(22, 6).Value = 1
$ws.Cells.Item(22, 7).Value = 0.1524443333333333
$ws.Cells.Item(22, 8).Value = 0.457333
$ws.Cells.Item(22, 9).Value = 0.0003829265417566354
$ws.Cells.Item(22, 10).Value = 0.0003829265417566354
$ws.Cells.Item(22, 13).Value = 2.655702666666667
$ws.Cells.Item(22, 14).Value = 7.967108
$ws.Cells.Item(22, 15).Value = 0.05887637219457465
$ws.Cells.Item(22, 16).Value = 0.05887637219457464
$ws.Cells.Item(22, 17).Value = 0.4048468225515555
$ws.Cells.Item(22, 18).Value = 3.643621402964
$ws.Cells.Item(22, 19).Value = 0.000022545325595645
$ws.Cells.Item(22, 20).Value = 0.000022545325595645
$ws.Cells.Item(23, 5).Value = 3
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(23, 7).Value = 0.1524443333333333
$ws.Cells.Item(23, 8).Value = 0.457333
$ws.Cells.Item(23, 9).Value = 0.0003829265417566354
$ws.Cells.Item(23, 10).Value = 0.0003829265417566354
$ws.Cells.Item(23, 15).Value = 0.111327724745791
$ws.Cells.Item(23, 16).Value = 0.111327724745791
$ws.Cells.Item(23, 17).Value = 0.7655138036745555
$ws.Cells.Item(23, 18).Value = 6.889624233070999
$ws.Cells.Item(23, 19).Value = 0.00004263034063854037
$ws.Cells.Item(23, 20).Value = 0.00004263034063854035
$ws.Cells.Item(24, 5).Value = 3
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = 0.1524443333333333
$ws.Cells.Item(24, 8).Value = 0.457333
$ws.Cells.Item(24, 9).Value = 0.0003829265417566354
$ws.Cells.Item(24, 10).Value = 0.0003829265417566354
$ws.Cells.Item(24, 13).Value = 0.9216453333333332
$ws.Cells.Item(24, 14).Value = 2.764936
$ws.Cells.Item(24, 15).Value = 0.02043268410948847
$ws.Cells.Item(24, 16).Value = 0.02043268410948846
$ws.Cells.Item(24, 17).Value = 0.1404996084097777
$ws.Cells.Item(24, 18).Value = 1.264496475688
$ws.Cells.Item(24, 19).Value = 0.000007824217064852177
$ws.Cells.Item(24, 20).Value = 0.000007824217064852175
$ws.Cells.Item(25, 5).Value = 3
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(25, 7).Value = 0.1524443333333333
$ws.Cells.Item(25, 8).Value = 0.457333
$ws.Cells.Item(25, 9).Value = 0.0003829265417566354
$ws.Cells.Item(25, 10).Value = 0.0003829265417566354
$ws.Cells.Item(25, 13).Value = 35.230657
$ws.Cells.Item(25, 14).Value = 105.691971
$ws.Cells.Item(25, 15).Value = 0.7810562907612387
$ws.Cells.Item(25, 16).Value = 0.7810562907612385
$ws.Cells.Item(25, 17).Value = 5.370714019260333
$ws.Cells.Item(25, 18).Value = 48.336426173343
$ws.Cells.Item(25, 19).Value = 0.0002990871843384663
$ws.Cells.Item(25, 20).Value = 0.0002990871843384662
$ws.Cells.Item(26, 5).Value = 3
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(26, 7).Value = 0.1524443333333333
$ws.Cells.Item(26, 8).Value = 0.457333
$ws.Cells.Item(26, 9).Value = 0.0003829265417566354
$ws.Cells.Item(26, 10).Value = 0.0003829265417566354
$ws.Cells.Item(26, 13).Value = 1.276824333333333
$ws.Cells.Item(26, 14).Value = 3.830473
$ws.Cells.Item(26, 15).Value = 0.02830692818890731
$ws.Cells.Item(26, 16).Value = 0.02830692818890731
$ws.Cells.Item(26, 17).Value = 0.1946446342787778
$ws.Cells.Item(26, 18).Value = 1.751801708509
$ws.Cells.Item(26, 19).Value = 0.0000108394741191317
$ws.Cells.Item(26, 20).Value = 0.00001083947411913169
